$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.301.84'
$ws.Range("E2").Value = '  -6.33%  '
$ws.Range("D3").Value = '2.468.18'
$ws.Range("E3").Value = '  -9.05%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '467.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.34'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.44%  '
$ws.Range("E8").Value = '  -5.94%  '
$ws.Range("D9").Value = '2.454.96'
$ws.Range("E9").Value = '  -9.65%  '
$ws.Range("E10").Value = '  -4.18%  '
$ws.Range("E11").Value = '  -8.53%  '
$ws.Range("E12").Value = '  -5.80%  '
$ws.Range("E13").Value = '  -3.37%  '
$ws.Range("D14").Value = '2.885.62'
$ws.Range("E14").Value = '  -9.82%  '
$ws.Range("D15").Value = '54.174.75'
$ws.Range("E15").Value = '  -6.82%  '
$ws.Range("E16").Value = '  +1.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.93'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.53%  '
$ws.Range("D18").Value = '2.457.48'
$ws.Range("E18").Value = '  -9.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -8.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '311.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.98%  '
$ws.Range("E21").Value = '  -11.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.72%  '
$ws.Range("E24").Value = '  -11.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '56.83'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.77%  '
$ws.Range("E26").Value = '  +0.61%  '
$ws.Range("E27").Value = '  -6.99%  '
$ws.Range("D28").Value = '2.541.01'
$ws.Range("E28").Value = '  -11.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.154'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.24'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("E32").Value = '  -6.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '150.22'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.88'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.15%  '
$ws.Range("E35").Value = '  -7.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.58'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -11.28%  '
$ws.Range("E38").Value = '  -2.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.809'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.85%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '33.72'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.97%  '
$ws.Range("E41").Value = '  -0.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.603'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.79%  '
$ws.Range("E43").Value = '  -2.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.08%  '
$ws.Range("E45").Value = '  -1.32%  '
$ws.Range("E46").Value = '  -4.85%  '
$ws.Range("D47").Value = '1.948.48'
$ws.Range("E47").Value = '  -8.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0219'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0876'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.76'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -8.96%  '
